$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.891504666666667
$ws.Range("H2").Value = 8.674514
$ws.Range("I2").Value = 0.1213590456377548
$ws.Range("J2").Value = 0.1213590456377548
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 325.332998030496
$ws.Range("R2").Value = 2927.996982274464
$ws.Range("S2").Value = 0.03974888570984019
$ws.Range("T2").Value = 0.03974888570984019

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.891504666666667
$ws.Range("H3").Value = 8.674514
$ws.Range("I3").Value = 0.1213590456377548
$ws.Range("J3").Value = 0.1213590456377548
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 307.4087745731746
$ws.Range("R3").Value = 2766.678971158572
$ws.Range("S3").Value = 0.03755892061574938
$ws.Range("T3").Value = 0.03755892061574938

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.891504666666667
$ws.Range("H4").Value = 8.674514
$ws.Range("I4").Value = 0.1213590456377548
$ws.Range("J4").Value = 0.1213590456377548
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 360.5465032907249
$ws.Range("R4").Value = 3244.918529616524
$ws.Range("S4").Value = 0.04405123931216519
$ws.Range("T4").Value = 0.0440512393121652

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.04042966666667
$ws.Range("H5").Value = 36.121289
$ws.Range("I5").Value = 0.505347637947847
$ws.Range("J5").Value = 0.505347637947847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 1354.709582934096
$ws.Range("R5").Value = 12192.38624640686
$ws.Range("S5").Value = 0.1655171676653134
$ws.Range("T5").Value = 0.1655171676653134

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.04042966666667
$ws.Range("H6").Value = 36.121289
$ws.Range("I6").Value = 0.505347637947847
$ws.Range("J6").Value = 0.505347637947847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 1280.071850422225
$ws.Range("R6").Value = 11520.64665380002
$ws.Range("S6").Value = 0.1563979983304588
$ws.Range("T6").Value = 0.1563979983304588

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.04042966666667
$ws.Range("H7").Value = 36.121289
$ws.Range("I7").Value = 0.505347637947847
$ws.Range("J7").Value = 0.505347637947847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 1501.341106061242
$ws.Range("R7").Value = 13512.06995455118
$ws.Range("S7").Value = 0.1834324719520748
$ws.Range("T7").Value = 0.1834324719520748

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.894099000000001
$ws.Range("H8").Value = 26.682297
$ws.Range("I8").Value = 0.3732933164143983
$ws.Range("J8").Value = 0.3732933164143982
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 1000.705247273808
$ws.Range("R8").Value = 9006.347225464271
$ws.Range("S8").Value = 0.1222652443617028
$ws.Range("T8").Value = 0.1222652443617028

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.894099000000001
$ws.Range("H9").Value = 26.682297
$ws.Range("I9").Value = 0.3732933164143983
$ws.Range("J9").Value = 0.3732933164143982
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 945.5713857361341
$ws.Range("R9").Value = 8510.142471625206
$ws.Range("S9").Value = 0.1155290399979581
$ws.Range("T9").Value = 0.1155290399979581

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.894099000000001
$ws.Range("H10").Value = 26.682297
$ws.Range("I10").Value = 0.3732933164143983
$ws.Range("J10").Value = 0.3732933164143982
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 1109.019927008545
$ws.Range("R10").Value = 9981.179343076903
$ws.Range("S10").Value = 0.1354990320547373
$ws.Range("T10").Value = 0.1354990320547373
